$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C13").Value = 11453
$ws.Range("C14").Value = 11197
$ws.Range("C15").Value = 11162
$ws.Range("C16:C19").Value = 11151
$ws.Range("C20:C31").Value = 11081
$ws.Range("C32:C39").Value = 10868
$ws.Range("C40").Value = 10660
$ws.Range("C41:C42").Value = 10254
$ws.Range("C43:C46").Value = 9919
$ws.Range("C47:C55").Value = 9481
$ws.Range("C56:C67").Value = 8942
$ws.Range("C68:C74").Value = 8575
$ws.Range("C75:C79").Value = 8397
$ws.Range("C80:C85").Value = 8095
$ws.Range("C86:C105").Value = 8018
$ws.Range("C106:C115").Value = 7968
$ws.Range("C116:C119").Value = 7900
$ws.Range("C120:C128").Value = 7748
$ws.Range("C129:C133").Value = 7295
